$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholders on every slide layout + the slide master:
#    8/6/2017 -> 8/11/2017
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($si = 1; $si -le $shapes.Count; $si++) {
        $sh = $shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "8/6/2017") {
                $tr.Text = "8/11/2017"
            }
        }
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

Update-DatePlaceholder $p.SlideMaster.Shapes

# ---------------------------------------------------------------------------
# Helper: replace a whole run's worth of text in-place (Find locates the
# exact same span as the existing run, so the assignment doesn't fragment
# the run into extra pieces).
# ---------------------------------------------------------------------------
function Replace-WholeRun($textRange, $oldText, $newText) {
    $found = $textRange.Find($oldText, 0)
    $found.Text = $newText
}

# ---------------------------------------------------------------------------
# 2) Slide 23 ("double-all" example): cons 11 -> cons 12, cons 22 -> cons 24
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$sh23 = $s23.Shapes.Item(2)
$tr23 = $sh23.TextFrame.TextRange

Replace-WholeRun $tr23 "(double-all (cons 11 empty)) " "(double-all (cons 12 empty)) "
Replace-WholeRun $tr23 "          = (cons 22 empty)" "          = (cons 24 empty)"
Replace-WholeRun $tr23 "(double-all (cons 33 (cons 11 empty)))" "(double-all (cons 33 (cons 12 empty)))"
Replace-WholeRun $tr23 "          = (cons 66 (cons 22 empty))" "          = (cons 66 (cons 24 empty))"

# ---------------------------------------------------------------------------
# 3) Slide 25 (step-by-step "double-all" trace): cons 11 -> cons 12,
#    cons 22 -> cons 24
# ---------------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$sh25 = $s25.Shapes.Item(2)
$tr25 = $sh25.TextFrame.TextRange

Replace-WholeRun $tr25 " (cons 11 (cons 22 (cons 33 empty))))" " (cons 12 (cons 22 (cons 33 empty))))"
Replace-WholeRun $tr25 "= (cons 22 (" "= (cons 24 ("
Replace-WholeRun $tr25 "= (cons 22 (cons 44 (" "= (cons 24 (cons 44 ("
Replace-WholeRun $tr25 "= (cons 22 (cons 44 (cons 66 (" "= (cons 24 (cons 44 (cons 66 ("

# Last line splits into 3 runs in the original deck: "= (" | "cons 22 " | "(cons 44 (cons 66 empty)))"
$lastFound = $tr25.Find("= (cons 22 (cons 44 (cons 66 empty)))", 0)
$lastSub = $tr25.Characters($lastFound.Start, $lastFound.Length)
$numFound = $lastSub.Find("cons 22 ", 0)
$numFound.Text = "cons 24 "

Write-Host "Edits applied"
